# Weekly price-sheet update ("Fruta / hortaliza, semanal"):
# A new week's record is prepended as row 159, pushing the former rows
# 159-232 down by two rows (to 161-234). The former row 160 retains its
# original J:Q values but gets the new week's date; the other rows just
# shift down keeping their data intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 159. Excel shifts
# rows 159:232 down to 161:234 and grows the used range automatically.
$ws.Rows("159:160").Insert()

# After the insert, the data that used to live in row 160 is now sitting
# in row 162. Copy that whole row back up into the (currently blank,
# except for the inherited date style) row 160 before we touch anything
# else, so every column besides the date is preserved verbatim.
$ws.Range("A162:R162").Copy()
$ws.Range("A160:R160").PasteSpecial()
$excel.CutCopyMode = 0

# Row 159: brand-new record for the new reporting week.
$ws.Cells.Item(159, 1).Value = 4
$ws.Cells.Item(159, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(159, 3).Value = "Los Lagos"
$ws.Cells.Item(159, 4).Value = 44572
$ws.Cells.Item(159, 5).Value = 10
$ws.Cells.Item(159, 6).Value = 100114014
$ws.Cells.Item(159, 7).Value = "Betarraga"
$ws.Cells.Item(159, 8).Value = "Sin especificar"
$ws.Cells.Item(159, 9).Value = "Primera"
$ws.Cells.Item(159, 10).Value = 110
$ws.Cells.Item(159, 11).Value = 10000
$ws.Cells.Item(159, 12).Value = 10000
$ws.Cells.Item(159, 13).Value = 10000
$ws.Cells.Item(159, 14).Value = "`$/malla 15 kilos"
$ws.Cells.Item(159, 15).Value = "Región Metropolitana"
$ws.Cells.Item(159, 16).Value = 667
$ws.Cells.Item(159, 17).Value = 15
$ws.Cells.Item(159, 18).Value = "Hortaliza"

# Row 160 only changes its reporting date; the rest came back via the
# copy/paste above.
$ws.Cells.Item(160, 4).Value = 44572

Write-Output "Inserted new weekly rows; used range is now $($ws.UsedRange.Address())"
